# Apply the StructureDefinition-quantity-with-conversion.xlsx update:
#  - bump Version 5.0.0 -> 6.0.0
#  - bump Date to the new publication timestamp
#  - replace the "Contact" metadata row with Publisher / Jurisdiction info
#  - remove the now-unused duplicate "Contact" row (net: one fewer row)
#  - clear the stray "N/A" RIM mapping value on the Elements sheet

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")

# Version: 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Date: 2021-12-16T17:36:56+00:00 -> 2022-01-21T20:46:54+00:00
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was blank; now "Alvearie Team"
$meta.Range("B9").Value = "Alvearie Team"

# Remove one of the duplicated "Contact" / "No display for ContactDetail" rows (row 10)
$meta.Rows.Item(10).Delete()

# The remaining row 10 (previously "Contact") becomes "Jurisdiction" / "United States of America"
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

$elements = $wb.Worksheets.Item("Elements")

# Quantity.extension (row 5) RIM Mapping (col AK) "N/A" -> blank
$elements.Range("AK5").Value = ""
